$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current newest-date header (currently in column B) before shifting.
$prevHeader = $ws.Range("B1").Value2

# Insert two new columns at C:D. This pushes the existing column C (and
# everything to its right) two columns to the right, i.e. old C -> new E.
# Column B (and A) are left untouched by the insert itself.
$ws.Range("C1:D1").EntireColumn.Insert()

# The two freshly inserted columns pick up column C's original custom width
# (8 characters wide) -- mirror that onto C, D and E just like Excel does
# when it clones formatting from the column being inserted next to.
$ws.Columns.Item(3).ColumnWidth = 7.15
$ws.Columns.Item(4).ColumnWidth = 7.15
$ws.Columns.Item(5).ColumnWidth = 7.15

# Column D now holds a blank cell where the old "Jun_13" header used to be
# (it got pushed there by the insert) -- restore it explicitly, then write
# the two new, more-recent snapshot headers into B1/C1.
$ws.Range("D1").Value = $prevHeader
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# For every data row, the two newly inserted cells (C, D) need the same
# placeholder rating ("UN") the rest of the sheet uses for "unchanged".
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
